$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2133.2896
$ws.Range("I15").Value = 2133.2896
$ws.Range("K15").Value = 6399.8688
$ws.Range("M15").Value = -6230.8688

$ws.Range("H40").Value = 1884.2858
$ws.Range("I40").Value = 1798
$ws.Range("J40").Value = 1932.2222
$ws.Range("K40").Value = 1798
$ws.Range("L40").Value = 1932.2222
$ws.Range("M40").Value = -1623
$ws.Range("N40").Value = -2282.2222

$ws.Range("H74").Value = 3800
$ws.Range("I74").Value = 3500
$ws.Range("J74").Value = 3842.8572
$ws.Range("K74").Value = 3500
$ws.Range("L74").Value = 3842.8572
$ws.Range("M74").Value = -2564
$ws.Range("N74").Value = -5714.8572

$ws.Range("H77").Value = 3800
$ws.Range("I77").Value = 3500
$ws.Range("J77").Value = 3842.8572
$ws.Range("K77").Value = 17500
$ws.Range("L77").Value = 19214.286
$ws.Range("M77").Value = -12820
$ws.Range("N77").Value = -28574.286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2464.63
$ws.Range("I32").Value = 2464.63
$ws.Range("K32").Value = 2464.63
$ws.Range("M32").Value = -2177.63

$ws.Range("H61").Value = 2344.718
$ws.Range("I61").Value = 1801.0526
$ws.Range("J61").Value = 2861.2
$ws.Range("K61").Value = 1801.0526
$ws.Range("L61").Value = 2861.2
$ws.Range("M61").Value = -1589.0526
$ws.Range("N61").Value = -3285.2

$ws.Range("H74").Value = 1706.6207
$ws.Range("I74").Value = 1179.9131
$ws.Range("J74").Value = 3725.6667
$ws.Range("K74").Value = 1179.9131
$ws.Range("L74").Value = 3725.6667
$ws.Range("M74").Value = -305.9131
$ws.Range("N74").Value = -5473.6667

$ws.Range("H77").Value = 1706.6207
$ws.Range("I77").Value = 1179.9131
$ws.Range("J77").Value = 3725.6667
$ws.Range("K77").Value = 5899.5655
$ws.Range("L77").Value = 18628.3335
$ws.Range("M77").Value = -1531.5655
$ws.Range("N77").Value = -27364.3335

$ws.Range("H132").Value = 1902.7164
$ws.Range("J132").Value = 5254.6665
$ws.Range("L132").Value = 15763.9995
$ws.Range("N132").Value = -20823.9995

$ws.Range("H134").Value = 30912.572
$ws.Range("J134").Value = 30912.572
$ws.Range("L134").Value = 30912.572
$ws.Range("N134").Value = -41052.572

$ws.Range("H136").Value = 2344.718
$ws.Range("I136").Value = 1801.0526
$ws.Range("J136").Value = 2861.2
$ws.Range("K136").Value = 5403.1578
$ws.Range("L136").Value = 8583.599999999999
$ws.Range("M136").Value = -2853.1578
$ws.Range("N136").Value = -13683.6

$ws.Range("H139").Value = 30481.666
$ws.Range("J139").Value = 30481.666
$ws.Range("L139").Value = 30481.666
$ws.Range("N139").Value = -40761.666

$ws.Range("H140").Value = 29656.666
$ws.Range("J140").Value = 29656.666
$ws.Range("L140").Value = 29656.666
$ws.Range("N140").Value = -40016.666

$ws.Range("H141").Value = 29410.75
$ws.Range("J141").Value = 29410.75
$ws.Range("L141").Value = 29410.75
$ws.Range("N141").Value = -39770.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 35000
$ws.Range("J9").Value = 35000
$ws.Range("L9").Value = 35000
$ws.Range("N9").Value = -35336

$ws.Range("H134").Value = 3482.366
$ws.Range("I134").Value = 3549.5518
$ws.Range("J134").Value = 3320
$ws.Range("K134").Value = 10648.6554
$ws.Range("L134").Value = 9960
$ws.Range("M134").Value = -8113.6554
$ws.Range("N134").Value = -15030

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9436524
$ws.Range("I58").Value = 1734.9
$ws.Range("J58").Value = 21742770
$ws.Range("K58").Value = 1734.9
$ws.Range("L58").Value = 21742770
$ws.Range("M58").Value = -1531.9
$ws.Range("N58").Value = -21743176

$ws.Range("H99").Value = 2645.182
$ws.Range("I99").Value = 1924.625
$ws.Range("J99").Value = 4566.6665
$ws.Range("K99").Value = 1924.625
$ws.Range("L99").Value = 4566.6665
$ws.Range("M99").Value = -426.625
$ws.Range("N99").Value = -7562.6665

$ws.Range("H126").Value = 2645.182
$ws.Range("I126").Value = 1924.625
$ws.Range("J126").Value = 4566.6665
$ws.Range("K126").Value = 5773.875
$ws.Range("L126").Value = 13699.9995
$ws.Range("M126").Value = -3303.875
$ws.Range("N126").Value = -18639.9995

$ws.Range("H132").Value = 4209.8335
$ws.Range("I132").Value = 3835.3333
$ws.Range("J132").Value = 5333.3335
$ws.Range("K132").Value = 11505.9999
$ws.Range("L132").Value = 16000.0005
$ws.Range("M132").Value = -8975.999899999999
$ws.Range("N132").Value = -21060.0005

$ws.Range("H134").Value = 16132026
$ws.Range("I134").Value = 26318342
$ws.Range("J134").Value = 3692.3333
$ws.Range("K134").Value = 78955026
$ws.Range("L134").Value = 11076.9999
$ws.Range("M134").Value = -78952491
$ws.Range("N134").Value = -16146.9999

$ws.Range("H136").Value = 9436524
$ws.Range("I136").Value = 1734.9
$ws.Range("J136").Value = 21742770
$ws.Range("K136").Value = 5204.700000000001
$ws.Range("L136").Value = 65228310
$ws.Range("M136").Value = -2654.700000000001
$ws.Range("N136").Value = -65233410

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3584.9656
$ws.Range("I94").Value = 1581
$ws.Range("J94").Value = 3905.6
$ws.Range("K94").Value = 4743
$ws.Range("L94").Value = 11716.8
$ws.Range("M94").Value = -4067
$ws.Range("N94").Value = -13068.8

$ws.Range("H131").Value = 1253
$ws.Range("J131").Value = 1035.7693
$ws.Range("L131").Value = 3107.3079
$ws.Range("N131").Value = -13187.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2601.4
$ws.Range("I7").Value = 2201.6
$ws.Range("J7").Value = 2801.3
$ws.Range("K7").Value = 2201.6
$ws.Range("L7").Value = 2801.3
$ws.Range("M7").Value = -2089.6
$ws.Range("N7").Value = -3025.3

$ws.Range("H40").Value = 15666.667
$ws.Range("I40").Value = 35000
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 35000
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -34864
$ws.Range("N40").Value = -6272

$ws.Range("H126").Value = 2601.4
$ws.Range("I126").Value = 2201.6
$ws.Range("J126").Value = 2801.3
$ws.Range("K126").Value = 6604.799999999999
$ws.Range("L126").Value = 8403.900000000001
$ws.Range("M126").Value = -4134.799999999999
$ws.Range("N126").Value = -13343.9

$ws.Range("H136").Value = 3676.3333
$ws.Range("I136").Value = 3781.8235
$ws.Range("J136").Value = 3538.3845
$ws.Range("K136").Value = 11345.4705
$ws.Range("L136").Value = 10615.1535
$ws.Range("M136").Value = -8795.470499999999
$ws.Range("N136").Value = -15715.1535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 50281.523
$ws.Range("I126").Value = 61006.59
$ws.Range("J126").Value = 4700
$ws.Range("K126").Value = 183019.77
$ws.Range("L126").Value = 14100
$ws.Range("M126").Value = -180549.77
$ws.Range("N126").Value = -19040

$ws.Range("H132").Value = 5941.2354
$ws.Range("I132").Value = 3067.238
$ws.Range("J132").Value = 10583.846
$ws.Range("K132").Value = 9201.714
$ws.Range("L132").Value = 31751.538
$ws.Range("M132").Value = -6671.714
$ws.Range("N132").Value = -36811.538

$ws.Range("H136").Value = 1944.3695
$ws.Range("I136").Value = 1786.24
$ws.Range("J136").Value = 2132.6191
$ws.Range("K136").Value = 5358.72
$ws.Range("L136").Value = 6397.8573
$ws.Range("M136").Value = -2808.72
$ws.Range("N136").Value = -11497.8573
